# Combine breastfeeding promotion delivery options
#   - Merge "Breastfeeding promotion - Health system" and
#     "Breastfeeding promotion - Home/community" into a single
#     "Breastfeeding promotion - Health system & Home/community" option.
#   - Add a "Complements group" header label on the
#     "OR stunting for complements" sheet.

$wb = $excel.ActiveWorkbook

# --- "OR stunting for complements": label the previously-blank header cell
$wsComplements = $wb.Worksheets.Item("OR stunting for complements")
$wsComplements.Range("A1").Value = "Complements group"

# --- "OR exclusiveBF by intervention": merge the two breastfeeding promotion
#     rows (row 2 = Health system, row 3 = Home/community) into row 2, then
#     remove the now-redundant row 3.
$wsExclusiveBF = $wb.Worksheets.Item("OR exclusiveBF by intervention")
$wsExclusiveBF.Range("A2").Value = "Breastfeeding promotion - Health system & Home/community"
$wsExclusiveBF.Rows(3).Delete()

# --- "Interventions coverages": merge the two breastfeeding promotion rows
#     (row 5 = Health system, row 6 = Home/community) into row 5, summing
#     their coverage values, then remove the now-redundant row 6.
$wsCoverages = $wb.Worksheets.Item("Interventions coverages")
$combinedCoverage = $wsCoverages.Range("B5").Value() + $wsCoverages.Range("B6").Value()
$wsCoverages.Range("A5").Value = "Breastfeeding promotion - Health system & Home/community"
$wsCoverages.Range("B5").Value = $combinedCoverage
$wsCoverages.Rows(6).Delete()
